$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1717705
$ws.Range("H2").Value = 0.343541
$ws.Range("I2").Value = 0.04063177891664595
$ws.Range("J2").Value = 0.02745976565347561
$ws.Range("M2").Value = 17.4294175
$ws.Range("N2").Value = 34.858835
$ws.Range("O2").Value = 0.4529581854295807
$ws.Range("P2").Value = 0.3776014560521451
$ws.Range("Q2").Value = 2.99385975868375
$ws.Range("R2").Value = 11.975439034735
$ws.Range("S2").Value = 0.01840449684885985
$ws.Range("T2").Value = 0.01036884749360308

# Row 3
$ws.Range("G3").Value = 0.1717705
$ws.Range("H3").Value = 0.343541
$ws.Range("I3").Value = 0.04063177891664595
$ws.Range("J3").Value = 0.02745976565347561
$ws.Range("O3").Value = 0.1017429801035258
$ws.Range("P3").Value = 0.127224648983019
$ws.Range("Q3").Value = 0.6724775567785001
$ws.Range("R3").Value = 4.034865340671
$ws.Range("S3").Value = 0.004133998273887167
$ws.Range("T3").Value = 0.003493559046419397

# Row 4
$ws.Range("G4").Value = 0.1717705
$ws.Range("H4").Value = 0.343541
$ws.Range("I4").Value = 0.04063177891664595
$ws.Range("J4").Value = 0.02745976565347561
$ws.Range("M4").Value = 5.397313
$ws.Range("N4").Value = 16.191939
$ws.Range("O4").Value = 0.1402661392829386
$ws.Range("P4").Value = 0.1753959862028526
$ws.Range("Q4").Value = 0.9270991526665
$ws.Range("R4").Value = 5.562594915999
$ws.Range("S4").Value = 0.005699262760835828
$ws.Range("T4").Value = 0.004816332677690574

# Row 5
$ws.Range("G5").Value = 0.1717705
$ws.Range("H5").Value = 0.343541
$ws.Range("I5").Value = 0.04063177891664595
$ws.Range("J5").Value = 0.02745976565347561
$ws.Range("M5").Value = 5.69137
$ws.Range("N5").Value = 11.38274
$ws.Range("O5").Value = 0.1479081344978025
$ws.Range("P5").Value = 0.1233012863987851
$ws.Range("Q5").Value = 0.9776094705849999
$ws.Range("R5").Value = 3.91043788234
$ws.Range("S5").Value = 0.006009770620888246
$ws.Range("T5").Value = 0.00338582442928272

# Row 6
$ws.Range("G6").Value = 0.1717705
$ws.Range("H6").Value = 0.343541
$ws.Range("I6").Value = 0.04063177891664595
$ws.Range("J6").Value = 0.02745976565347561
$ws.Range("M6").Value = 1.988496
$ws.Range("N6").Value = 5.965488000000001
$ws.Range("O6").Value = 0.05167731737988258
$ws.Range("P6").Value = 0.06461997237892773
$ws.Range("Q6").Value = 0.341564952168
$ws.Range("R6").Value = 2.049389713008
$ws.Range("S6").Value = 0.002099741334784735
$ws.Range("T6").Value = 0.001774449298059423

# Row 7
$ws.Range("G7").Value = 0.1717705
$ws.Range("H7").Value = 0.343541
$ws.Range("I7").Value = 0.04063177891664595
$ws.Range("J7").Value = 0.02745976565347561
$ws.Range("M7").Value = 4.057513666666667
$ws.Range("N7").Value = 12.172541
$ws.Range("O7").Value = 0.1054472433062699
$ws.Range("P7").Value = 0.1318566499842704
$ws.Range("Q7").Value = 0.6969611512801667
$ws.Range("R7").Value = 4.181766907681
$ws.Range("S7").Value = 0.004284509077390133
$ws.Range("T7").Value = 0.003620752708420425

# Row 8
$ws.Range("I8").Value = 0.9593682210833541
$ws.Range("J8").Value = 0.9725402343465244
$ws.Range("M8").Value = 17.4294175
$ws.Range("N8").Value = 34.858835
$ws.Range("O8").Value = 0.4529581854295807
$ws.Range("P8").Value = 0.3776014560521451
$ws.Range("Q8").Value = 70.6888545725175
$ws.Range("R8").Value = 424.133127435105
$ws.Range("S8").Value = 0.4345536885807209
$ws.Range("T8").Value = 0.3672326085585421

# Row 9
$ws.Range("I9").Value = 0.9593682210833541
$ws.Range("J9").Value = 0.9725402343465244
$ws.Range("O9").Value = 0.1017429801035258
$ws.Range("P9").Value = 0.127224648983019
$ws.Range("S9").Value = 0.09760898182963861
$ws.Range("T9").Value = 0.1237310899365996

# Row 10
$ws.Range("I10").Value = 0.9593682210833541
$ws.Range("J10").Value = 0.9725402343465244
$ws.Range("M10").Value = 5.397313
$ws.Range("N10").Value = 16.191939
$ws.Range("O10").Value = 0.1402661392829386
$ws.Range("P10").Value = 0.1753959862028526
$ws.Range("Q10").Value = 21.889995677673
$ws.Range("R10").Value = 197.009961099057
$ws.Range("S10").Value = 0.1345668765221028
$ws.Range("T10").Value = 0.170579653525162

# Row 11
$ws.Range("I11").Value = 0.9593682210833541
$ws.Range("J11").Value = 0.9725402343465244
$ws.Range("M11").Value = 5.69137
$ws.Range("N11").Value = 11.38274
$ws.Range("O11").Value = 0.1479081344978025
$ws.Range("P11").Value = 0.1233012863987851
$ws.Range("Q11").Value = 23.08260882777
$ws.Range("R11").Value = 138.49565296662
$ws.Range("S11").Value = 0.1418983638769143
$ws.Range("T11").Value = 0.1199154619695024

# Row 12
$ws.Range("I12").Value = 0.9593682210833541
$ws.Range("J12").Value = 0.9725402343465244
$ws.Range("M12").Value = 1.988496
$ws.Range("N12").Value = 5.965488000000001
$ws.Range("O12").Value = 0.05167731737988258
$ws.Range("P12").Value = 0.06461997237892773
$ws.Range("Q12").Value = 8.064784985616001
$ws.Range("R12").Value = 72.58306487054401
$ws.Range("S12").Value = 0.04957757604509785
$ws.Range("T12").Value = 0.06284552308086831

# Row 13
$ws.Range("I13").Value = 0.9593682210833541
$ws.Range("J13").Value = 0.9725402343465244
$ws.Range("M13").Value = 4.057513666666667
$ws.Range("N13").Value = 12.172541
$ws.Range("O13").Value = 0.1054472433062699
$ws.Range("P13").Value = 0.1318566499842704
$ws.Range("Q13").Value = 16.456143385687
$ws.Range("R13").Value = 148.105290471183
$ws.Range("S13").Value = 0.1011627342288798
$ws.Range("T13").Value = 0.12823589727585
